$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.898.39"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "'1.641.17"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.61%  "

$ws.Range("D5").Value = "'215.55"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("D6").Value = "'0.5034"
$ws.Range("E6").Value = "  -2.15%  "

$ws.Range("E7").Value = "  -0.74%  "

$ws.Range("E8").Value = "  -0.97%  "

$ws.Range("D9").Value = "'0.06382"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("D10").Value = "'19.53"
$ws.Range("E10").Value = "  -1.79%  "

$ws.Range("D11").Value = "'0.07744"
$ws.Range("E11").Value = "  -1.30%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.650.81"
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.262"
$ws.Range("E13").Value = "  -1.00%  "

$ws.Range("D14").Value = "'1.864.58"
$ws.Range("E14").Value = "  -1.15%  "

$ws.Range("D15").Value = "'0.5457"
$ws.Range("E15").Value = "  -1.34%  "

$ws.Range("D16").Value = "'0.0₅7905"
$ws.Range("E16").Value = "  -1.42%  "

$ws.Range("D17").Value = "'64.28"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").Value = "'25.900.05"
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  -0.77%  "

$ws.Range("D20").Value = "'202.92"
$ws.Range("E20").Value = "  -3.33%  "

$ws.Range("D21").Value = "'4.397"
$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("D22").Value = "'9.880"
$ws.Range("E22").Value = "  -2.36%  "

$ws.Range("D23").Value = "'5.971"
$ws.Range("E23").Value = "  -1.28%  "

$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  -0.72%  "

$ws.Range("D25").Value = "'1.868"
$ws.Range("E25").Value = "  +3.19%  "

$ws.Range("D26").Value = "'140.68"
$ws.Range("E26").Value = "  -2.88%  "

$ws.Range("D27").Value = "'0.1134"
$ws.Range("E27").Value = "  -3.53%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'15.66"
$ws.Range("E28").Value = "  -1.36%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'6.788"
$ws.Range("E29").Value = "  -3.07%  "

$ws.Range("D30").Value = "'1.244"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").Value = "'0.04980"
$ws.Range("E31").Value = "  -2.16%  "

$ws.Range("D32").Value = "'3.276"
$ws.Range("E32").Value = "  -2.54%  "

$ws.Range("D33").Value = "'3.196"
$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").Value = "'1.542"
$ws.Range("E34").Value = "  -1.19%  "

$ws.Range("D35").Value = "'2.366"
$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("D36").Value = "'2.632"
$ws.Range("E36").Value = "  -4.08%  "

$ws.Range("D37").Value = "'0.8916"
$ws.Range("E37").Value = "  -3.58%  "

$ws.Range("D38").Value = "'1.153.80"
$ws.Range("E38").Value = "  -1.64%  "

$ws.Range("D39").Value = "'0.5605"
$ws.Range("E39").Value = "  -1.86%  "

$ws.Range("D40").Value = "'0.01563"
$ws.Range("E40").Value = "  -1.76%  "

$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("D42").Value = "'5.682"
$ws.Range("E42").Value = "  +0.28%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'99.85"
$ws.Range("E43").Value = "  -0.57%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8073"
$ws.Range("E44").Value = "  -2.63%  "

$ws.Range("D45").Value = "'1.777.00"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("E46").Value = "  +4.98%  "

$ws.Range("D47").Value = "'0.4543"
$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D48").Value = "'1.006"
$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").Value = "'54.78"
$ws.Range("E49").Value = "  -1.40%  "

$ws.Range("D50").Value = "'0.05056"
$ws.Range("E50").Value = "  -0.61%  "

$ws.Range("D51").Value = "'0.9991"
$ws.Range("E51").Value = "  -0.87%  "
